$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) cells we are about to rewrite to stay as plain
# text (matching the original inlineStr storage) instead of being
# auto-coerced to numbers by Excel when the new value parses as a number.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated Price (D) and Volume(1h) (E) values row by row.
$ws.Range("D2").Value = '30.127.00'
$ws.Range("E2").Value = '  +0.45%  '
$ws.Range("D3").Value = '1.918.38'
$ws.Range("E3").Value = '  +2.67%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '320.12'
$ws.Range("E5").Value = '  +0.34%  '
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").Value = '0.5057'
$ws.Range("E7").Value = '  -0.65%  '
$ws.Range("D8").Value = '0.4077'
$ws.Range("E8").Value = '  +3.96%  '
$ws.Range("D9").Value = '0.08353'
$ws.Range("E9").Value = '  +2.41%  '
$ws.Range("D10").Value = '42.43'
$ws.Range("E10").Value = '  +0.95%  '
$ws.Range("E11").Value = '  +1.96%  '
$ws.Range("D12").Value = '24.07'
$ws.Range("E12").Value = '  +5.94%  '
$ws.Range("D13").Value = '6.413'
$ws.Range("E13").Value = '  +2.81%  '
$ws.Range("D14").Value = '1.901.57'
$ws.Range("E14").Value = '  +1.94%  '
$ws.Range("D15").Value = '7.247'
$ws.Range("E15").Value = '  +1.42%  '
$ws.Range("D16").Value = '0.9993'
$ws.Range("E16").Value = '  -0.19%  '
$ws.Range("D17").Value = '92.54'
$ws.Range("E17").Value = '  +1.09%  '
$ws.Range("D18").Value = '0.00001094'
$ws.Range("E18").Value = '  +1.61%  '
$ws.Range("D19").Value = '0.06509'
$ws.Range("E19").Value = '  +2.31%  '
$ws.Range("D20").Value = '18.54'
$ws.Range("E20").Value = '  +3.91%  '
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("D22").Value = '5.945'
$ws.Range("E22").Value = '  +2.80%  '
$ws.Range("D23").Value = '30.146.56'
$ws.Range("E23").Value = '  +0.58%  '
$ws.Range("D24").Value = '11.36'
$ws.Range("E24").Value = '  +2.72%  '
$ws.Range("D25").Value = '2.192'
$ws.Range("E25").Value = '  +0.24%  '
$ws.Range("D26").Value = '2.128.80'
$ws.Range("E26").Value = '  +2.28%  '
$ws.Range("D27").Value = '21.84'
$ws.Range("D28").Value = '162.90'
$ws.Range("E28").Value = '  +1.09%  '
$ws.Range("D29").Value = '2.282'
$ws.Range("E29").Value = '  +3.00%  '
$ws.Range("D30").Value = '128.76'
$ws.Range("E30").Value = '  +1.58%  '
$ws.Range("D31").Value = '1.144'
$ws.Range("E31").Value = '  +9.20%  '
$ws.Range("D32").Value = '0.1045'
$ws.Range("E32").Value = '  +1.36%  '
$ws.Range("D33").Value = '5.951'
$ws.Range("E33").Value = '  +1.20%  '
$ws.Range("D34").Value = '3.780'
$ws.Range("E34").Value = '  +1.32%  '
$ws.Range("D35").Value = '0.02458'
$ws.Range("E35").Value = '  +2.07%  '
$ws.Range("D36").Value = '5.366'
$ws.Range("E36").Value = '  +3.33%  '
$ws.Range("D37").Value = '0.06447'
$ws.Range("E37").Value = '  +2.22%  '
$ws.Range("E38").Value = '  +1.08%  '
$ws.Range("D39").Value = '0.6538'
$ws.Range("E39").Value = '  +4.15%  '
$ws.Range("E40").Value = '  +2.31%  '
$ws.Range("D41").Value = '8.627'
$ws.Range("E41").Value = '  +1.93%  '
$ws.Range("D42").Value = '11.41'
$ws.Range("E42").Value = '  +1.52%  '
$ws.Range("D43").Value = '1.212'
$ws.Range("E43").Value = '  +0.67%  '
$ws.Range("D44").Value = '13.46'
$ws.Range("E44").Value = '  +4.50%  '
$ws.Range("D45").Value = '0.6081'
$ws.Range("E45").Value = '  +3.54%  '
$ws.Range("D46").Value = '2.191'
$ws.Range("E46").Value = '  +10.70%  '
$ws.Range("D47").Value = '3.624'
$ws.Range("E47").Value = '  +0.09%  '
$ws.Range("D48").Value = '1.210'
$ws.Range("E48").Value = '  +0.91%  '
$ws.Range("D49").Value = '122.21'
$ws.Range("E49").Value = '  -0.28%  '
$ws.Range("D50").Value = '79.10'
$ws.Range("E50").Value = '  +3.73%  '
$ws.Range("D51").Value = '1.138'
$ws.Range("E51").Value = '  -0.71%  '
